# Nexial command-reference workbook update:
# - add new `image` command: ocr(image,saveVar)
# - add new `tn.5250` command category/group with its 5259/5250-terminal commands
# - adjust the `image` and `target` named ranges to account for the new rows
# - insert a new column for the `tn.5250` command list ahead of the `web` column
#   (shifting web/webalert/webcookie/ws/ws.async/xml one column to the right)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Make room for the new `tn.5250` column by inserting a blank column at Z
#    (column 26). Everything from Z..AE (web..xml) shifts right to AA..AF.
# ---------------------------------------------------------------------------
$ws.Columns.Item(26).Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new Z column: header + the tn.5250 command list.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = "tn.5250"
$ws.Range("Z2").Value = "close(profile)"
$ws.Range("Z3").Value = "open(profile)"
$ws.Range("Z4").Value = "saveText(profile,var)"
$ws.Range("Z5").Value = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value = "updateScreenFields(profile)"

# ---------------------------------------------------------------------------
# 3) Update the `image` (column K) command list:
#    - K2 changes from colorbit(source,bit,saveTo) to colorbit(image,bit,saveTo)
#    - a new ocr(image,saveVar) row is inserted before resize(...)/saveDiff(...),
#      pushing those two down by one row.
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = "colorbit(image,bit,saveTo)"
$ws.Range("K8").Value = $ws.Range("K7").Value2
$ws.Range("K7").Value = $ws.Range("K6").Value2
$ws.Range("K6").Value = "ocr(image,saveVar)"

# ---------------------------------------------------------------------------
# 4) Update the `target` (column A) category list: insert "tn.5250" in its
#    alphabetically-sorted position (between "step" and "web"), pushing
#    web/webalert/webcookie/ws/ws.async/xml down by one row.
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = $ws.Range("A31").Value
$ws.Range("A31").Value = $ws.Range("A30").Value
$ws.Range("A30").Value = $ws.Range("A29").Value
$ws.Range("A29").Value = $ws.Range("A28").Value
$ws.Range("A28").Value = $ws.Range("A27").Value
$ws.Range("A27").Value = $ws.Range("A26").Value
$ws.Range("A26").Value = "tn.5250"

# ---------------------------------------------------------------------------
# 5) Fix up the defined names so they continue to point at the right ranges.
# ---------------------------------------------------------------------------
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
